# Generate Report for Handoff
# Bumps the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# forward for the rows that were re-generated, and marks their Priority as "ht".

$wb = $excel.ActiveWorkbook

$rows = @(8, 10, 11, 12, 13, 14)

# --- Overview sheet: "Latest HO Xliff Generate Date" column G ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-08-13 12:26:20"
}

# --- zh-cn / de-de sheets: "Latest Handoff Datetime" column H, "Priority" column E ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($r in $rows) {
        $ws.Range("H$r").Value = "2016-08-13 12:26:12"
        $ws.Range("E$r").Value = "ht"
    }
}
